$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.76"
$ws.Range("E2").Value = "'2.55%"
$ws.Range("G2").Value = "'4"

$ws.Range("D3").Value = "'39.52"
$ws.Range("E3").Value = "'-0.73%"
$ws.Range("G3").Value = "'4"

$ws.Range("D4").Value = "'5.133"
$ws.Range("E4").Value = "'-0.29%"
$ws.Range("G4").Value = "'4"

$ws.Range("D5").Value = "'0.08178"
$ws.Range("E5").Value = "'0.79%"
$ws.Range("G5").Value = "'4"

$ws.Range("D6").Value = "'1.989"
$ws.Range("E6").Value = "'2.17%"
$ws.Range("G6").Value = "'4"

$ws.Range("B7").Value = "'KuCoinToken"
$ws.Range("C7").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.309"
$ws.Range("E7").Value = "'1.90%"
$ws.Range("G7").Value = "'4"

$ws.Range("B8").Value = "'MXToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9376"
$ws.Range("E8").Value = "'0.94%"
$ws.Range("G8").Value = "'4"

$ws.Range("B9").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1304"
$ws.Range("E9").Value = "'-8.68%"
$ws.Range("G9").Value = "'4"

$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1970"
$ws.Range("E10").Value = "'2.35%"
$ws.Range("G10").Value = "'4"

$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09099"
$ws.Range("E11").Value = "'-0.35%"
$ws.Range("G11").Value = "'4"

$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03477"
$ws.Range("E12").Value = "'-0.88%"
$ws.Range("G12").Value = "'4"

$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09757"
$ws.Range("E13").Value = "'-0.68%"
$ws.Range("G13").Value = "'4"

$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001406"
$ws.Range("E14").Value = "'0.72%"
$ws.Range("G14").Value = "'4"

$ws.Range("B15").Value = "'TigerCash"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006171"
$ws.Range("E15").Value = "'3.72%"
$ws.Range("G15").Value = "'4"

$ws.Range("B16").Value = "'LEO"
$ws.Range("C16").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.635"
$ws.Range("E16").Value = "'-7.69%"
$ws.Range("G16").Value = "'4"

$ws.Range("B17").Value = "'GateToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.366"
$ws.Range("E17").Value = "'3.28%"
$ws.Range("G17").Value = "'4"

$ws.Range("D18").Value = "'3.300"
$ws.Range("E18").Value = "'-2.27%"
$ws.Range("G18").Value = "'4"

$ws.Range("E19").Value = "'1.86%"
$ws.Range("G19").Value = "'4"

$ws.Range("D20").Value = "'0.1315"
$ws.Range("E20").Value = "'-0.37%"
$ws.Range("G20").Value = "'4"

$ws.Range("E21").Value = "'6.55%"
$ws.Range("G21").Value = "'4"

$ws.Range("D22").Value = "'0.2578"
$ws.Range("E22").Value = "'6.46%"
$ws.Range("G22").Value = "'4"

$ws.Range("D23").Value = "'0.04359"
$ws.Range("E23").Value = "'-0.49%"
$ws.Range("G23").Value = "'4"

$ws.Range("E24").Value = "'1.00%"
$ws.Range("G24").Value = "'4"

$ws.Range("D25").Value = "'0.004771"
$ws.Range("E25").Value = "'9.55%"
$ws.Range("G25").Value = "'4"

$ws.Range("D26").Value = "'0.0003889"
$ws.Range("E26").Value = "'198.84%"
$ws.Range("G26").Value = "'4"

$ws.Range("E27").Value = "'-7.54%"
$ws.Range("G27").Value = "'4"

$ws.Range("G28").Value = "'4"

$ws.Range("G29").Value = "'4"

$ws.Range("G30").Value = "'4"

$ws.Range("G31").Value = "'4"

$ws.Range("G32").Value = "'4"

$ws.Range("G33").Value = "'4"

$ws.Range("G34").Value = "'4"

$ws.Range("G35").Value = "'4"

$ws.Range("G36").Value = "'4"

$ws.Range("G37").Value = "'4"

$ws.Range("G38").Value = "'4"

$ws.Range("D39").Value = "'0.02203"
$ws.Range("E39").Value = "'8.08%"
$ws.Range("G39").Value = "'4"

$ws.Range("D40").Value = "'0.05193"
$ws.Range("E40").Value = "'2.50%"
$ws.Range("G40").Value = "'4"

$ws.Range("D41").Value = "'0.007766"
$ws.Range("E41").Value = "'5.01%"
$ws.Range("G41").Value = "'4"

$ws.Range("D42").Value = "'0.01036"
$ws.Range("E42").Value = "'6.31%"
$ws.Range("G42").Value = "'4"

$ws.Range("G43").Value = "'4"

$ws.Range("D44").Value = "'0.002099"
$ws.Range("E44").Value = "'-1.54%"
$ws.Range("G44").Value = "'4"

$ws.Range("D45").Value = "'0.009749"
$ws.Range("E45").Value = "'2.71%"
$ws.Range("G45").Value = "'4"

$ws.Range("D46").Value = "'0.00006726"
$ws.Range("E46").Value = "'6.05%"
$ws.Range("G46").Value = "'4"

$ws.Range("E47").Value = "'-0.03%"
$ws.Range("G47").Value = "'4"

$ws.Range("D48").Value = "'0.002882"
$ws.Range("E48").Value = "'5.29%"
$ws.Range("G48").Value = "'4"

$ws.Range("D49").Value = "'0.001689"
$ws.Range("E49").Value = "'29.91%"
$ws.Range("G49").Value = "'4"

$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("G50").Value = "'4"

$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.03%"
$ws.Range("G51").Value = "'4"
